# Insert a new data row at position 422 (shifting existing rows 422-524 down to 423-525)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 422 (this shifts rows down and keeps formatting)
$ws.Rows.Item(422).Insert()

# Populate the newly inserted row 422 with the new record
$ws.Cells.Item(422, 1).Value = 3
$ws.Cells.Item(422, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(422, 3).Value = "Coquimbo"
$ws.Cells.Item(422, 4).Value = (Get-Date -Year 2023 -Month 1 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(422, 5).Value = 5
$ws.Cells.Item(422, 6).Value = 100112017
$ws.Cells.Item(422, 7).Value = "Apio"
$ws.Cells.Item(422, 8).Value = "Americana (o)"
$ws.Cells.Item(422, 9).Value = "Primera"
$ws.Cells.Item(422, 10).Value = 110
$ws.Cells.Item(422, 11).Value = 12000
$ws.Cells.Item(422, 12).Value = 12000
$ws.Cells.Item(422, 13).Value = 12000
$ws.Cells.Item(422, 14).Value = "`$/docena de matas"
$ws.Cells.Item(422, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(422, 16).Value = 2000
$ws.Cells.Item(422, 17).Value = 6
$ws.Cells.Item(422, 18).Value = "Hortaliza"
